# Update the "Command Strings List" sheet to reflect changes to the
# drivetrain's "Both Side" command row (row 4):
#   - Value column (E4):      "00 to 99 + 00 to 99"  -> "left,right"
#   - Notes column (F4):      old left/right-wheel note -> new ±100 note
#   - Example column (G4):    "D1575"                -> "D-40,40"
# The old, now-orphaned note that used to live in E3
# ("Two Digit Positive Number 00 to 99") is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the stale note in E3 (cell becomes blank / removed from the sheet).
$ws.Range("E3").ClearContents()

# Row 4 ("Both Side" drivetrain command) gets the new value/notes/example.
$ws.Range("E4").Value = "left,right"
$ws.Range("F4").Value = "negative is backwards, and positive is forwards. Range is " + [char]0x00B1 + "100"
$ws.Range("G4").Value = '"D-40,40"'

# Match the author's active selection after editing.
$ws.Range("E4").Select()
